$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace the sample data row with header-like text values
$ws.Range("A2").Value = "SerialBird"
$ws.Range("B2").Value = "Strain"
$ws.Range("C2").Value = "Sub Species"
$ws.Range("D2").Value = "Date of Bird"
$ws.Range("E2").Value = "Gender"
$ws.Range("F2").Value = "Cage  Number"
$ws.Range("G2").Value = "Father"
$ws.Range("H2").Value = "Mother"

# Row 3: replace sample data row with numeric 100 values
$ws.Range("A3").Value = 100
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 100

# Row 4: remove entirely (delete the row so the used range shrinks to A1:H3)
$ws.Range("A4:H4").Delete()
